$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new header cells for team record
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match formatting of the other header cells (bold, bordered, centered)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill data rows 2-48 with the team's record: 72 wins, 90 losses, 0 ties
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 72  # AD
    $ws.Cells.Item($r, 31).Value = 90  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
